# Auto-generated PowerShell Excel COM-interop script
# Applies odds/value updates to row 3 and row 4 of Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 1.25
$ws.Range("H3").Value = 6.25
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 1.67
$ws.Range("K3").Value = 2.75
$ws.Range("L3").Value = 8.5
$ws.Range("N3").Value = 19
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 5.5
$ws.Range("Q3").Value = 1.44
$ws.Range("R3").Value = 2.7
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.75
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.83
$ws.Range("W3").Value = 9
$ws.Range("Y3").Value = 9
$ws.Range("Z3").Value = 8
$ws.Range("AB3").Value = 23
$ws.Range("AC3").Value = 19
$ws.Range("AD3").Value = 12
$ws.Range("AE3").Value = 23
$ws.Range("AF3").Value = 51
$ws.Range("AG3").Value = 251
$ws.Range("AH3").Value = 26
$ws.Range("AI3").Value = 51
$ws.Range("AJ3").Value = 26
$ws.Range("AK3").Value = 101
$ws.Range("AL3").Value = 51
$ws.Range("AM3").Value = 51
$ws.Range("AN3").Value = 3.4
$ws.Range("AO3").Value = 5.5
$ws.Range("AQ3").Value = 13
$ws.Range("AT3").Value = 3.75
$ws.Range("AU3").Value = 9.5
$ws.Range("AX3").Value = 10
$ws.Range("AY3").Value = 41
$ws.Range("BA3").Value = 151
$ws.Range("BB3").Value = 151
$ws.Range("BC3").Value = 251
# Row 4
$ws.Range("J4").Value = 2.6
$ws.Range("K4").Value = 2.25
$ws.Range("M4").Value = 1.04
$ws.Range("N4").Value = 13
$ws.Range("O4").Value = 1.22
$ws.Range("P4").Value = 4
$ws.Range("Q4").Value = 1.8
$ws.Range("R4").Value = 2
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 10
$ws.Range("AB4").Value = 23
$ws.Range("AC4").Value = 12
$ws.Range("AE4").Value = 13
$ws.Range("AG4").Value = 151
$ws.Range("AP4").Value = 19
$ws.Range("AS4").Value = 126
$ws.Range("AU4").Value = 7.5
$ws.Range("AY4").Value = 19
$ws.Range("BA4").Value = 51
